$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: intervention_type, styled like the other header cells (bold/centered/bordered)
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "intervention_type"

# Data rows
$ws.Range("K2").Value = "OTHER"
$ws.Range("K3").Value = "DEVICE"
$ws.Range("K4").Value = "DRUG"
$ws.Range("K5").Value = "BIOLOGICAL"
$ws.Range("K6").Value = "DRUG"
